# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# Add new player "T.Smith" row to the RB sheet (row 7) with zeroed stats
$rbSheet = $wb.Worksheets.Item("RB")
$rbSheet.Range("A7").Value = "T.Smith"
$rbSheet.Range("B7:J7").Value = 0

# Make RB the active sheet/tab (simulating the user navigating there), and
# set its selection to J8 as the last-used cell.
$rbSheet.Activate()
$rbSheet.Range("J8").Select()

# The WR sheet is no longer the active tab; make sure its selection stays as before.
$wrSheet = $wb.Worksheets.Item("WR")
$wrSheet.Range("J11").Select()

# Re-activate RB so it is the tab shown/selected when the workbook is saved.
$rbSheet.Activate()
